$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 16:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1689581
$ws.Range("C4").Value = 3145
$ws.Range("E4").Value = 1138455
$ws.Range("G4").Value = 81
$ws.Range("H4").Value = 99381

# Row 13 - India
$ws.Range("B13").Value = 141228
$ws.Range("C13").Value = 2692
$ws.Range("D13").Value = 58727
$ws.Range("E13").Value = 78444
$ws.Range("G13").Value = 33
$ws.Range("H13").Value = 4057

# Row 31 - Suiza
$ws.Range("E31").Value = 733
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = 1913

# Row 40 - Rumania
$ws.Range("E40").Value = 5451
$ws.Range("G40").Value = 17
$ws.Range("H40").Value = 1202

# Rows 62-64: Moldavia moved up ahead of Australia and Armenia in the ranking.
# Row 62 now shows Moldavia's refreshed data, row 63 now shows the former
# Australia row data, row 64 now shows the former Armenia row data.
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 7147
$ws.Range("C62").Value = 54
$ws.Range("D62").Value = 3802
$ws.Range("E62").Value = 3084
$ws.Range("G62").Value = 11
$ws.Range("H62").Value = 261

$ws.Range("A63").Value = "Australia"
$ws.Range("B63").Value = 7118
$ws.Range("C63").Value = 4
$ws.Range("D63").Value = 6531
$ws.Range("E63").Value = 485
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 102

$ws.Range("A64").Value = "Armenia"
$ws.Range("B64").Value = 7113
$ws.Range("C64").Value = 452
$ws.Range("D64").Value = 3145
$ws.Range("E64").Value = 3881
$ws.Range("H64").Value = 87

# Row 153 - Yemen
$ws.Range("B153").Value = 233
$ws.Range("C153").Value = 11
$ws.Range("E153").Value = 179
$ws.Range("G153").Value = 2
$ws.Range("H153").Value = 44
